# Add team record (Wins / Losses / Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labels in AD1:AF1 ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header style (bold, centered, bordered) by copying
# the format from the adjacent header cell (AC1) onto the new header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows (2-41): constant team record for every player row ---
$lastRow = 41
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("AD$r").Value = 90
    $ws.Range("AE$r").Value = 72
    $ws.Range("AF$r").Value = 0
}
